$d = $word.ActiveDocument

# 1) Portuguese "Programa" list: "1. ...;2. ...;...;7. ...;"
#    Insert a manual line break before every "N." that directly follows a
#    semicolon (keeping the semicolon with the preceding sentence).
$rng1 = $d.Paragraphs(14).Range
$rng1.Find.Execute(";([0-9]\.)", $true, $false, $true, $false, $false, $true, 1, $false, ";^l\1", 2)

# 2) English "Programa" list: "1. ...2. ...3. ...7. ..." (no separating
#    punctuation). Insert a manual line break between a lowercase letter
#    and the following "N. ".
$rng2 = $d.Paragraphs(15).Range
$rng2.Find.Execute("([a-z])([0-9]\. )", $true, $false, $true, $false, $false, $true, 1, $false, "\1^l\2", 2)

# 3) Bibliografia list: "1. ...New Jersey.2. ...GmbH.3. ...USA."
#    Insert a manual line break before every "N." that directly follows a
#    period (keeping the period with the preceding sentence).
$rng3 = $d.Paragraphs(19).Range
$rng3.Find.Execute("\.([0-9]\.)", $true, $false, $true, $false, $false, $true, 1, $false, ".^l\1", 2)
